$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Run ID (corrected translation)
$ws.Range("B2").Value = "e91c8e84-e7d7-471d-83bd-f6628178a777_36"

# Updated metric values
$ws.Range("B3").Value = 0.85714000000000001
$ws.Range("B4").Value = 0.7
$ws.Range("B5").Value = 0.86121999999999999
$ws.Range("B6").Value = 0.7
$ws.Range("B7").Value = 0.71436999999999995
$ws.Range("B8").Value = 0.83762000000000003
$ws.Range("B9").Value = 0.81747999999999998
$ws.Range("B10").Value = 0.72499999999999998
$ws.Range("B11").Value = 0.69771000000000005
$ws.Range("B12").Value = 0.85714000000000001
$ws.Range("B13").Value = 0.82838999999999996
$ws.Range("B14").Value = 0.51561000000000001
$ws.Range("B15").Value = 0.45954
$ws.Range("B16").Value = 0.45
$ws.Range("B17").Value = 0.70643
$ws.Range("B18").Value = 0.85714000000000001
$ws.Range("B19").Value = 0.83816000000000002
$ws.Range("B20").Value = 0.72499999999999998
$ws.Range("B21").Value = 0.85714000000000001
$ws.Range("B22").Value = 0.85714000000000001
$ws.Range("B23").Value = 0.89319000000000004

# Update the sheet's current selection (was C4 -> now E6)
$ws.Range("E6").Select()
